$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values we still need before we start overwriting cells,
# since column A/B will be overwritten with data currently living in B/C.
$row1_newB = $ws.Cells.Item(1, 3).Value2   # old C1 -> new B1

$labels = @{}
$firstVals = @{}
for ($r = 2; $r -le 6; $r++) {
    $labels[$r]    = $ws.Cells.Item($r, 2).Value2   # old column B (label)
    $firstVals[$r] = $ws.Cells.Item($r, 3).Value2   # old column C (first data value)
}

# Clear everything from C1:J6 - removes values AND styling, and shrinks
# the sheet's used range/dimension down automatically.
$ws.Range("C1:J6").Clear()

# Row 1: only B1 survives, carrying the old C1 value. Style s="1" is
# already present on B1 so no style change needed.
$ws.Cells.Item(1, 2).Value = $row1_newB

# Rows 2-6: column A becomes the old label (inherits old column A's bold
# style, already s="1"), column B becomes the old first data value
# (no style, like old column B/C).
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $labels[$r]
    $ws.Cells.Item($r, 2).Value = $firstVals[$r]
}

# Row 4's B cell must stay the literal text "21.22" rather than being
# reinterpreted as a number. Use the apostrophe text-prefix trick, then
# reset the style back to Normal so no extra numFmt/quotePrefix style
# lingers on the cell.
$ws.Cells.Item(4, 2).Value = "'21.22"
$ws.Cells.Item(4, 2).Style = "Normal"
